# csv_format_standard_demonstration_20200825.xlsx
# "small updates, added crosswalk"
#
# The Standard/Description table in row 6 is reworded for ten of the
# twelve "Standard Element" columns (Character Set, File Name, Data
# Matrix, Variable Name, Units, Missing Data Value, Flags, Temporal
# Data, Timestamps, Spatial Data). The Delimiter (D6) and Consistent
# Values (H6) descriptions are unchanged. Finally the active selection
# in the sheet is moved from E6 down to A31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 'Data stored using the comma separated values format (CSV) (RFC 4180) must use the the standard seven-bit American Standard Code for Information Interchange (US-ASCII) characters (RFC 20).'
$ws.Range("C6").Value = 'Unique file names must be used. No spaces. Only letters, numbers, a hyphen "-" (ASCII Code 45), and an underscore "_" (ASCII Code 95) can be used in file names.'
$ws.Range("E6").Value = 'The contents of the file must be organized in a logical and readable matrix format. There can be no empty lines or rows in the file, and the file must contain the same number of columns across all of its rows.'
$ws.Range("F6").Value = 'Unique variable names must be used. No spaces. Letters, numbers, a hyphen "-" (ASCII Code 45), and an underscore "_" (ASCII Code 95) are preferred in variable names.'
$ws.Range("G6").Value = 'Provide the units of measurement for the variable in the variable name following the same naming conventions for the variable. If units are not provided here, it must be documented elsewhere. Data should be represented with units of measurement approved by the International System of Units (SP 330), derived units (e.g., degree Celsius), or non-SI units accepted for use with SI (e.g., mixing ratio). Explanations of units the do not conform to the international standards must be documented elsewhere.'
$ws.Range("I6").Value = 'If a cell does not contain a value, a missing data value must be indicated. Missing data must be represented by values that can never be construed as actual data and must be consistent across variables. For columns containing numeric data, "-9999" is preferred as the missing data value or use the correct precision given the data in the column. For columns containing character data, the string "NA" is preferred as the missing data value. Explanations for individual missing values can be reported as a separate variable (i.e., in an adjacent column). If a coding system is used to describe the missing data value, it must be documented elsewhere.'
$ws.Range("J6").Value = 'Measurement uncertainty, limits of detection, data quality indicators, and other flags pertaining to individual values should be reported as a separate variable (i.e., in an adjacent column) but only in addition to the reported values. If a coding system is used to describe the flags, it must be documented elsewhere.'
$ws.Range("K6").Value = 'All dates and times must be reported in Coordinated Universal Time (UTC) and follow the ISO 8601 standard (RFC 3339). Note that the use of "Z" and "T" characters are unnecessary. All times must be preceded with a date. In cases where the entire file consists of temporal data collected at a single date and time, the date and time must be documented elsewhere if not provided as a variable. Temporal data using different standards can be provided as a separate variable (i.e., in an adjacent column) but only in addition to UTC format.'
$ws.Range("L6").Value = 'For data with multiple timestamped records or when applicable, the variable name should specify if the measurement is the start, stop, or midpoint value, or it should be documented elsewhere.'
$ws.Range("M6").Value = 'All geographic coordinates must be provided in WGS84 decimal format (EPSG 4326). Latitude and longitude must be provided as separate variables (i.e., in an adjacent column). For geolocated records, each row must contain coordinates. In cases where the entire file consists of measurements collected at a single location, a pair of geographic coordinates must be documented elsewhere if not provided as variables. Spatial data using different standards can be provided as a separate variable (i.e., in an adjacent column) but only in addition to WGS84 decimal format.'

# Move the selection/scroll position to match the saved view state (A31,
# scrolled so row 7 is at the top).
$ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
